$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.773.30"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").Value = "1.634.91"
$ws.Range("E3").Value = "  -0.03%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.66"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("E6").Value = "  -0.88%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0634"
$ws.Range("E9").Value = "  -1.23%  "

$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("E12").Value = "  +0.01%  "

$ws.Range("D13").Value = "1.861.24"
$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("D14").Value = "1.639.11"
$ws.Range("E14").Value = "  +0.21%  "

$ws.Range("E15").Value = "  +0.62%  "

$ws.Range("D16").Value = "0.0₃0763"
$ws.Range("E16").Value = "  -0.59%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.21"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").Value = "25.804.51"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.46"
$ws.Range("E20").Value = "  +1.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.46"
$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("E22").Value = "  +0.49%  "

$ws.Range("E23").Value = "  +2.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.84"
$ws.Range("E24").Value = "  +4.40%  "

$ws.Range("E25").Value = "  -0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.49"
$ws.Range("E26").Value = "  +1.42%  "

$ws.Range("E27").Value = "  +1.83%  "

$ws.Range("E28").Value = "  +0.97%  "

$ws.Range("E29").Value = "  +0.27%  "

$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("E34").Value = "  -0.65%  "

$ws.Range("E35").Value = "  -0.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.904"
$ws.Range("E36").Value = "  +0.18%  "

$ws.Range("D37").Value = "1.130.23"
$ws.Range("E37").Value = "  +1.15%  "

$ws.Range("E38").Value = "  -1.86%  "

$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("E40").Value = "  -1.10%  "

$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("E42").Value = "  +0.89%  "

$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.75"
$ws.Range("E44").Value = "  +1.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.801"
$ws.Range("E45").Value = "  +0.23%  "

$ws.Range("D46").Value = "1.770.04"
$ws.Range("E46").Value = "  -0.20%  "

$ws.Range("E47").Value = "  +2.24%  "

$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("E49").Value = "  -0.94%  "

$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("E51").Value = "  +4.09%  "

